# QUBES Code First commit
# Adds two new data rows (5 and 6) to the "TestData" sheet, each with a
# single value in column C, extending the used range from A1:J4 to A1:J6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$ws.Range("C5").Value = "Var1-VS1P320220133"
$ws.Range("C6").Value = "PRK033103"
